# Update the cryptocurrency price/volume table with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps storing values as text, matching the
# original "inline string" cell type used for figures like "30.532.04".
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.532.04"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.911.86"
$ws.Range("E3").Value = "  -1.57%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "239.48"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  -1.69%  "

$ws.Range("E8").Value = "  -2.54%  "

$ws.Range("D9").Value = "0.06693"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("D10").Value = "18.69"
$ws.Range("E10").Value = "  -3.68%  "

$ws.Range("D11").Value = "101.22"
$ws.Range("E11").Value = "  -3.67%  "

$ws.Range("D12").Value = "0.07685"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("D13").Value = "1.912.01"
$ws.Range("E13").Value = "  -1.66%  "

$ws.Range("D14").Value = "5.214"
$ws.Range("E14").Value = "  -1.70%  "

$ws.Range("D15").Value = "0.6691"
$ws.Range("E15").Value = "  -3.76%  "

$ws.Range("D16").Value = "30.542.80"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").Value = "256.30"
$ws.Range("E17").Value = "  -6.61%  "

$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").Value = "0.000007470"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("D20").Value = "12.65"
$ws.Range("E20").Value = "  -3.55%  "

$ws.Range("D21").Value = "5.382"
$ws.Range("E21").Value = "  -0.92%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").Value = "6.295"
$ws.Range("E23").Value = "  -2.32%  "

$ws.Range("D24").Value = "9.335"
$ws.Range("E24").Value = "  -3.67%  "

$ws.Range("D25").Value = "167.00"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").Value = "19.15"
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("D27").Value = "2.059"
$ws.Range("E27").Value = "  -4.80%  "

$ws.Range("D28").Value = "4.764"
$ws.Range("E28").Value = "  +5.12%  "

$ws.Range("D29").Value = "0.1008"
$ws.Range("E29").Value = "  -2.80%  "

$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").Value = "4.244"
$ws.Range("E32").Value = "  -2.43%  "

$ws.Range("D33").Value = "0.04716"
$ws.Range("E33").Value = "  -2.46%  "

$ws.Range("D34").Value = "0.7295"
$ws.Range("E34").Value = "  -2.07%  "

$ws.Range("E35").Value = "  -3.78%  "

$ws.Range("D36").Value = "0.9998"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -0.86%  "

$ws.Range("D38").Value = "0.01915"
$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("D39").Value = "2.614"
$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("D40").Value = "6.240"
$ws.Range("E40").Value = "  -2.92%  "

$ws.Range("D41").Value = "74.49"
$ws.Range("E41").Value = "  -3.44%  "

$ws.Range("E42").Value = "  -5.17%  "

$ws.Range("D43").Value = "0.8619"
$ws.Range("E43").Value = "  -3.56%  "

$ws.Range("D44").Value = "105.47"
$ws.Range("E44").Value = "  -2.28%  "

$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "0.4234"
$ws.Range("E46").Value = "  -3.74%  "

$ws.Range("D47").Value = "7.363"
$ws.Range("E47").Value = "  -4.95%  "

$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("D50").Value = "908.04"
$ws.Range("E50").Value = "  -9.22%  "

$ws.Range("D51").Value = "8.739"
$ws.Range("E51").Value = "  -4.13%  "
